# This script applies the cryptocurrency price/volume refresh described
# by the commit "Updated cryptos list ... with GitHub Actions".
# Rows 21-24 also got re-sorted (Uniswap/ShibaInu and PancakeSwap/Litecoin swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '42.906.04'
$ws.Range("E2").Value = '  -5.26%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.215.73'
$ws.Range("E3").Value = '  -6.47%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '319.11'
$ws.Range("E5").Value = '  +0.14%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.23'
$ws.Range("E6").Value = '  -9.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.578'
$ws.Range("E7").Value = '  -8.92%  '

$ws.Range("E8").Value = '  -0.10%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.554'
$ws.Range("E9").Value = '  -9.75%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.54'
$ws.Range("E10").Value = '  -10.84%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.88'
$ws.Range("E11").Value = '  -3.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0823'
$ws.Range("E12").Value = '  -10.45%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.58'
$ws.Range("E13").Value = '  -10.62%  '

$ws.Range("E14").Value = '  -1.83%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.554.74'
$ws.Range("E15").Value = '  -6.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.858'
$ws.Range("E16").Value = '  -12.56%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.29'
$ws.Range("E17").Value = '  -7.44%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.221.34'
$ws.Range("E18").Value = '  -6.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '42.839.08'
$ws.Range("E19").Value = '  -5.31%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.00'
$ws.Range("E20").Value = '  -9.79%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("B21").Value = 'ShibaInu'
$ws.Range("C21").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D21").Value = '0.0₃0959'
$ws.Range("E21").Value = '  -9.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '6.50'
$ws.Range("E22").Value = '  -10.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("B23").Value = 'PancakeSwap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D23").Value = '3.17'
$ws.Range("E23").Value = '  -12.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '64.68'
$ws.Range("E24").Value = '  -11.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '234.83'
$ws.Range("E25").Value = '  -11.30%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.14'
$ws.Range("E26").Value = '  -8.47%  '

$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.02'
$ws.Range("E28").Value = '  +0.90%  '

$ws.Range("E29").Value = '  -2.97%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.91'
$ws.Range("E30").Value = '  -11.85%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.28'
$ws.Range("E31").Value = '  -16.25%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '35.43'
$ws.Range("E32").Value = '  -5.12%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.12'
$ws.Range("E33").Value = '  -10.11%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0859'
$ws.Range("E34").Value = '  -9.40%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '151.44'
$ws.Range("E35").Value = '  -10.33%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.65'
$ws.Range("E36").Value = '  -8.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.14'
$ws.Range("E37").Value = '  +1.89%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.121'
$ws.Range("E38").Value = '  -8.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.89'
$ws.Range("E39").Value = '  -1.59%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.38'
$ws.Range("E40").Value = '  -7.07%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.103'
$ws.Range("E41").Value = '  -11.70%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.64'
$ws.Range("E42").Value = '  -9.87%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0318'
$ws.Range("E43").Value = '  -10.07%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.50'
$ws.Range("E44").Value = '  +4.08%  '

$ws.Range("E45").Value = '  +0.04%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.726.97'
$ws.Range("E46").Value = '  -7.57%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.202'
$ws.Range("E47").Value = '  -11.43%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '83.43'
$ws.Range("E48").Value = '  -15.88%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.25'
$ws.Range("E49").Value = '  -12.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.76'
$ws.Range("E50").Value = '  -4.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '73.64'
$ws.Range("E51").Value = '  -12.90%  '
